$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2024-03-21 Thursday" "2024-03-22 Friday"

Replace-Text "583×3=" "204×3="
Replace-Text "599×4=" "509×5="
Replace-Text "300×7=" "268×4="
Replace-Text "594×9=" "537×4="
Replace-Text "290×4=" "993×2="
Replace-Text "945×3=" "492×8="
Replace-Text "115×7=" "767×5="
Replace-Text "781×5=" "424×4="
Replace-Text "116×9=" "797×6="
Replace-Text "697×6=" "441×5="
Replace-Text "733×4=" "607×2="
Replace-Text "440×7=" "265×3="
Replace-Text "367×2=" "720×4="
Replace-Text "354×3=" "473×8="
Replace-Text "213×5=" "299×4="
Replace-Text "931×7=" "843×4="
Replace-Text "813×4=" "148×3="
Replace-Text "239×8=" "194×7="
Replace-Text "918×5=" "393×8="
Replace-Text "708×4=" "492×3="
Replace-Text "395×7=" "930×7="
Replace-Text "115×6=" "692×9="
Replace-Text "410×8=" "348×6="
Replace-Text "333×4=" "388×9="
Replace-Text "937×7=" "668×4="

Write-Output "done"
